$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, pushing existing rows 13..58 down to 14..59.
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the new weekly data entry.
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(13, 3).Value = "Bíobío"
$ws.Cells.Item(13, 4).Value = "2021-12-31"
$ws.Cells.Item(13, 5).Value = 8
$ws.Cells.Item(13, 6).Value = 100112001
$ws.Cells.Item(13, 7).Value = "Berenjena"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 250
$ws.Cells.Item(13, 11).Value = 8000
$ws.Cells.Item(13, 12).Value = 9000
$ws.Cells.Item(13, 13).Value = 8400
$ws.Cells.Item(13, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(13, 15).Value = "Región Metropolitana"
$ws.Cells.Item(13, 16).Value = 140
$ws.Cells.Item(13, 17).Value = 60
$ws.Cells.Item(13, 18).Value = "Hortaliza"
